$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Articulos")
$ws2 = $wb.Worksheets.Item("Catalogo")

# Fix the text in Catalogo!A5 to remove the accent (Artículo -> Articulo)
$ws2.Range("A5").Value = "Articulo"

# Hide the Catalogo sheet
$ws2.Visible = $false

# Update selections to match target state
$ws2.Range("A5").Select()
$ws1.Range("A2").Select()
